{"js": "// Update the date line and the 26 multiplication problems in the practice\n// table to the next day's generated worksheet content.\nconst pairs = [\n  [\"2024-11-27 Wednesday\", \"2024-11-28 Thursday\"],\n\n  [\"753\u00d78=\", \"477\u00d72=\"],\n  [\"202\u00d74=\", \"310\u00d77=\"],\n  [\"845\u00d78=\", \"360\u00d73=\"],\n  [\"339\u00d78=\", \"845\u00d74=\"],\n  [\"402\u00d72=\", \"700\u00d75=\"],\n\n  [\"445\u00d78=\", \"895\u00d74=\"],\n  [\"881\u00d74=\", \"242\u00d74=\"],\n  [\"136\u00d75=\", \"125\u00d74=\"],\n  [\"676\u00d79=\", \"182\u00d77=\"],\n  [\"681\u00d76=\", \"398\u00d72=\"],\n\n  [\"336\u00d74=\", \"675\u00d76=\"],\n  [\"529\u00d74=\", \"346\u00d79=\"],\n  [\"903\u00d77=\", \"332\u00d72=\"],\n  [\"602\u00d76=\", \"632\u00d78=\"],\n  [\"768\u00d75=\", \"395\u00d79=\"],\n\n  [\"796\u00d77=\", \"637\u00d73=\"],\n  [\"878\u00d75=\", \"310\u00d77=\"],\n  [\"871\u00d76=\", \"817\u00d78=\"],\n  [\"542\u00d73=\", \"239\u00d77=\"],\n  [\"675\u00d78=\", \"920\u00d74=\"],\n\n  [\"466\u00d78=\", \"218\u00d76=\"],\n  [\"213\u00d73=\", \"326\u00d74=\"],\n  [\"657\u00d73=\", \"214\u00d78=\"],\n  [\"896\u00d76=\", \"932\u00d74=\"],\n  [\"719\u00d79=\", \"494\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 26 multiplication problems in the practice\n# table to the next day's generated worksheet content.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $range = $d.Content\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n\nReplace-Text \"2024-11-27 Wednesday\" \"2024-11-28 Thursday\"\n\nReplace-Text \"753\u00d78=\" \"477\u00d72=\"\nReplace-Text \"202\u00d74=\" \"310\u00d77=\"\nReplace-Text \"845\u00d78=\" \"360\u00d73=\"\nReplace-Text \"339\u00d78=\" \"845\u00d74=\"\nReplace-Text \"402\u00d72=\" \"700\u00d75=\"\n\nReplace-Text \"445\u00d78=\" \"895\u00d74=\"\nReplace-Text \"881\u00d74=\" \"242\u00d74=\"\nReplace-Text \"136\u00d75=\" \"125\u00d74=\"\nReplace-Text \"676\u00d79=\" \"182\u00d77=\"\nReplace-Text \"681\u00d76=\" \"398\u00d72=\"\n\nReplace-Text \"336\u00d74=\" \"675\u00d76=\"\nReplace-Text \"529\u00d74=\" \"346\u00d79=\"\nReplace-Text \"903\u00d77=\" \"332\u00d72=\"\nReplace-Text \"602\u00d76=\" \"632\u00d78=\"\nReplace-Text \"768\u00d75=\" \"395\u00d79=\"\n\nReplace-Text \"796\u00d77=\" \"637\u00d73=\"\nReplace-Text \"878\u00d75=\" \"310\u00d77=\"\nReplace-Text \"871\u00d76=\" \"817\u00d78=\"\nReplace-Text \"542\u00d73=\" \"239\u00d77=\"\nReplace-Text \"675\u00d78=\" \"920\u00d74=\"\n\nReplace-Text \"466\u00d78=\" \"218\u00d76=\"\nReplace-Text \"213\u00d73=\" \"326\u00d74=\"\nReplace-Text \"657\u00d73=\" \"214\u00d78=\"\nReplace-Text \"896\u00d76=\" \"932\u00d74=\"\nReplace-Text \"719\u00d79=\" \"494\u00d73=\"\n"}
